$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Record V2")

# --- Bagging V1 block (rows 59-66): update date + accuracy results ---
$ws.Range("A59").Value = 45955
$ws.Range("D59").Value = 0.7273
$ws.Range("A60").ClearContents()
$ws.Range("D60").Value = 0.7273
$ws.Range("D61").Value = 0.8182
$ws.Range("D62").Value = 0.7273
$ws.Range("D63").Value = 0.7273
$ws.Range("D64").Value = 0.8182
$ws.Range("D65").Value = 0.8182
$ws.Range("D66").Value = 0.7273

# --- Clear the stray date / run-note cells in column A for rows 67-74 ---
# (cutting from the A column of rows 75-82, which have no cell there at all,
#  reproduces the exact "vanish" / "blank stub" pattern Excel leaves behind)
$ws.Range("A75:A82").Cut($ws.Range("A67"))

# --- Fix the mislabelled "Bagging V2" block (rows 67-74): it is actually
#     the Random Forest results block ---
$ws.Range("B67").Value = "Random "
$ws.Range("B67").Style = $ws.Range("B75").Style
$ws.Range("B68").Value = "Forest"
$ws.Range("B68").Style = $ws.Range("B76").Style
$ws.Range("D67").ClearContents()
$ws.Range("D68").ClearContents()

# --- Relabel rows 75-76: this used to be "Random Forest" which just moved
#     up to 67-68; rows 75-76 are actually the "Boosting" results block ---
$ws.Range("B75").Value = "Boosting "
$ws.Range("B75").Style = $ws.Range("B83").Style
$ws.Range("B76").ClearContents()
$ws.Range("B76").Style = $ws.Range("B84").Style

# --- Remove the now fully-duplicated trailing "Boosting" rows ---
$ws.Rows("83:90").Delete()

# --- Update sheet view state ---
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("E63").Select()
